# "Generate Report for Archive"
#
# The localization-status report was regenerated: the two open status
# records moved from "Ready for handoff" to "In Translation", and the
# Status column (mirrored on the Overview sheet as the per-locale
# "zh-cn"/"de-de" columns, and as column C on each per-locale sheet)
# narrows accordingly to fit the new, shorter text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "In Translation" -------------
$overview.Range("E2:F3").Value = "In Translation"
$zhcn.Range("C2:C3").Value     = "In Translation"
$dede.Range("C2:C3").Value     = "In Translation"

# --- Narrow the Status columns to match the shorter text ---------------
# Excel's Range.ColumnWidth setter quantizes to whole screen pixels
# (width is internally tracked as (pixels-5)/6 "characters"), so the
# nearest reachable column width to the report's target of
# 13.4101845877511 characters is obtained by feeding in the character
# width whose pixel-rounded result lands closest to it.
$targetWidth = 13.4101845877511
$nearestStep = [Math]::Round($targetWidth * 6 - 5)
$inputWidth = $nearestStep / 6

$overview.Columns.Item(5).ColumnWidth = $inputWidth
$overview.Columns.Item(6).ColumnWidth = $inputWidth
$zhcn.Columns.Item(3).ColumnWidth = $inputWidth
$dede.Columns.Item(3).ColumnWidth = $inputWidth
